$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Update the input assumption in B11 (discount/cap rate parameter) from 0.15 to 0.125.
# This cell feeds the formulas in B20:B120 (and the mirrored formulas on the
# CSC-CSCCCMvSoECBtY sheet + chart caches), so Excel will recalc everything downstream.
$ws.Range("B11").Value = 0.125

# Update the saved view state for the About sheet: scroll so row 4 is the
# top-left visible row, and select cell B12 instead of the old A4:XFD4 selection.
$ws.Activate()
$ws.Range("B12").Select()
$excel.ActiveWindow.ScrollRow = 4

$wb.Save()
